$d = $word.ActiveDocument

$replacements = @(
    @{old = "461×9="; new = "618×6="},
    @{old = "244×2="; new = "643×9="},
    @{old = "517×3="; new = "670×8="},
    @{old = "153×5="; new = "440×5="},
    @{old = "547×2="; new = "778×3="},
    @{old = "682×4="; new = "486×7="},
    @{old = "526×5="; new = "390×3="},
    @{old = "474×5="; new = "541×8="},
    @{old = "238×2="; new = "266×6="},
    @{old = "257×3="; new = "779×5="},
    @{old = "433×4="; new = "990×3="},
    @{old = "103×9="; new = "170×4="},
    @{old = "362×9="; new = "740×8="},
    @{old = "808×3="; new = "634×9="},
    @{old = "592×3="; new = "162×7="},
    @{old = "725×7="; new = "243×5="},
    @{old = "483×5="; new = "249×8="},
    @{old = "128×5="; new = "471×3="},
    @{old = "366×5="; new = "755×5="},
    @{old = "310×4="; new = "966×6="},
    @{old = "929×3="; new = "722×9="},
    @{old = "423×3="; new = "695×8="},
    @{old = "414×9="; new = "481×2="},
    @{old = "328×8="; new = "829×5="},
    @{old = "926×2="; new = "692×3="}
)

foreach ($r in $replacements) {
    $d.Content.Find.Execute($r.old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $r.new, 2)
}
